# Insert a new weekly price record as row 322 in the "Perejil" sheet,
# pushing the existing rows 322-365 down to 323-366 (dimension grows
# from A1:R365 to A1:R366).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(322).Insert()

$ws.Cells.Item(322, 1).Value()  = 6
$ws.Cells.Item(322, 2).Value()  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(322, 3).Value()  = "Metropolitana"
$ws.Cells.Item(322, 4).Value()  = 44491
$ws.Cells.Item(322, 5).Value()  = 13
$ws.Cells.Item(322, 6).Value()  = 100112044
$ws.Cells.Item(322, 7).Value()  = "Perejil"
$ws.Cells.Item(322, 8).Value()  = "Sin especificar"
$ws.Cells.Item(322, 9).Value()  = "Primera"
$ws.Cells.Item(322, 10).Value() = 270
$ws.Cells.Item(322, 11).Value() = 7000
$ws.Cells.Item(322, 12).Value() = 8000
$ws.Cells.Item(322, 13).Value() = 7444
$ws.Cells.Item(322, 14).Value() = "$/docena de atados"
$ws.Cells.Item(322, 15).Value() = "Región Metropolitana"
$ws.Cells.Item(322, 16).Value() = 2481
$ws.Cells.Item(322, 17).Value() = 3
$ws.Cells.Item(322, 18).Value() = "Hortaliza"
